# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
$range.Value2 = 45188
